$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values after repulling data / recalculating means
$updates = @{
    3  = 1
    6  = -1
    7  = 3
    8  = 0
    9  = 2
    14 = -2
    17 = -6
    18 = -2
    27 = -3
    28 = -3
    30 = -2
    33 = 1
    37 = 2
    39 = 3
    46 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
